# Apply the "S2026" dated edit + new Overview section described by the diff.
$d = $word.ActiveDocument

# 1. "F2025" -> "S2026" in the Date paragraph.
$d.Content.Find.Execute("F2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "S2026", 2) | Out-Null

# 2. "Debate stuff" -> "Nothing to see here yet" in the FirstParagraph paragraph.
$d.Content.Find.Execute("Debate stuff", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Nothing to see here yet", 2) | Out-Null

# 3. Insert a new "Overview" Heading2 paragraph right after the Date paragraph
#    (i.e. directly before the FirstParagraph paragraph that now reads
#    "Nothing to see here yet").
$datePara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Style.NameLocal -eq "Date") {
        $datePara = $para
    }
}

$datePara.Range.InsertParagraphAfter() | Out-Null
$overviewPara = $datePara.Next()
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$overviewPara.Range.InsertXML("<w:p $ns><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t xml:space='preserve'>Overview</w:t></w:r></w:p>") | Out-Null

# Re-locate the paragraphs after the InsertXML reshuffle, then wrap the new
# "Overview" heading paragraph and the following "Nothing to see here yet"
# paragraph in a bookmark named "overview".
$overviewPara = $null
$bodyPara = $null
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text.TrimEnd([char]13)
    if ($overviewPara -eq $null -and $para.Style.NameLocal -eq "Heading 2" -and $text -eq "Overview") {
        $overviewPara = $para
    } elseif ($overviewPara -ne $null -and $bodyPara -eq $null) {
        $bodyPara = $para
    }
}

$bmRange = $d.Range($overviewPara.Range.Start, $bodyPara.Range.End)
$d.Bookmarks.Add("overview", $bmRange) | Out-Null
